# Append one new telemetry row to each of the four worksheets, matching
# the format of the existing rows (date-time in col A, hex-byte text in
# columns B-E, and numeric counters in columns F-I).
#
# NOTE: this runtime's PowerShell-like dialect does not bind named
# (-paramName value) arguments on user-defined functions, so the helper
# below uses plain positional parameters.

$wb = $excel.ActiveWorkbook

function Add-DataRow($ws, $row, $aVal, $bVal, $cVal, $dVal, $eVal, $fVal, $gVal, $hVal, $iVal) {
    # Column A: date/time value, formatted like the rest of the column.
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 1).Value = $aVal

    # Columns B-E: plain text (hex byte strings).
    $ws.Cells.Item($row, 2).Value = $bVal
    $ws.Cells.Item($row, 3).Value = $cVal
    $ws.Cells.Item($row, 4).Value = $dVal
    $ws.Cells.Item($row, 5).Value = $eVal

    # Column F: plain numeric.
    $ws.Cells.Item($row, 6).Value = $fVal

    # Column G: numeric, unless it cannot round-trip through a double
    # cleanly, in which case it is kept as text (mirrors source data).
    if ($gVal -is [string]) {
        $ws.Cells.Item($row, 7).NumberFormat = "@"
        $ws.Cells.Item($row, 7).Value = $gVal
        $ws.Cells.Item($row, 7).Style = "Normal"
    } else {
        $ws.Cells.Item($row, 7).Value = $gVal
    }

    # Columns H-I: plain numeric.
    $ws.Cells.Item($row, 8).Value = $hVal
    $ws.Cells.Item($row, 9).Value = $iVal
}

# --- Sheet 1: ROW50-FE-LIFTER -> new row 75 ---
$ws1 = $wb.Worksheets.Item(1)
Add-DataRow $ws1 75 45761.75865222222 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x4e" "0xe" 400 568631262647114000000000.0 334 14

# --- Sheet 2: ROW50-MID-LIFTER -> new row 77 ---
$ws2 = $wb.Worksheets.Item(2)
Add-DataRow $ws2 77 45761.72050925926 "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x52" "0x19" 400 "568631262647113771663628" 338 25

# --- Sheet 3: ROW11-FE-LIFTER -> new row 75 ---
$ws3 = $wb.Worksheets.Item(3)
Add-DataRow $ws3 75 45761.78976460648 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x01,0x4e" "0x14" 400 568631262647114000000000.0 334 20

# --- Sheet 4: ROW11-MID-LIFTER -> new row 75 ---
$ws4 = $wb.Worksheets.Item(4)
Add-DataRow $ws4 75 45761.91668688657 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x01,0x56" "0x19" 400 568631262647114000000000.0 342 25

Write-Host "Rows appended to all four sheets."
